$wb = $excel.ActiveWorkbook

# Workbook-level view change: tabRatio 963 -> 925 (tabRatio is stored in the
# OOXML as parts-per-thousand, exposed via COM as a 0..1 ratio)
$wb.Windows.Item(1).TabRatio = 0.925

# Select the "parsed mile posts" sheet (sheetId=1, first sheet) and update selection
$ws = $wb.Worksheets.Item("parsed mile posts")
$ws.Activate()
$ws.Range("C6").Select()

# Update data values in row 2
$ws.Range("B2").Value = 9.96
$ws.Range("C2").Value = 10.56
$ws.Range("D2").Value = 147000
$ws.Range("F2").Value = 2
$ws.Range("G2").Value = 2

# Remove the Comments cell H2 (was s="1" t="n" v="5") entirely
$ws.Range("H2").Clear()
